$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster data (player, position, team) for rows 2-18
$data = @(
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Nick Richards", "C", "Phoenix Suns"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
